$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8030.5386
$ws.Range("J51").Value = 8981.727999999999
$ws.Range("L51").Value = 8981.727999999999
$ws.Range("N51").Value = -9949.727999999999

$ws.Range("H113").Value = 9353.058999999999
$ws.Range("I113").Value = 6665.8335
$ws.Range("K113").Value = 6665.8335
$ws.Range("M113").Value = -3411.8335

$ws.Range("H129").Value = 1524.6471
$ws.Range("I129").Value = 953.9091
$ws.Range("J129").Value = 2571
$ws.Range("K129").Value = 2861.7273
$ws.Range("L129").Value = 7713
$ws.Range("M129").Value = 2138.2727
$ws.Range("N129").Value = -17713

$ws.Range("H138").Value = 3186.5293
$ws.Range("I138").Value = 2442.3
$ws.Range("J138").Value = 3666.6775
$ws.Range("K138").Value = 7326.900000000001
$ws.Range("L138").Value = 11000.0325
$ws.Range("M138").Value = -2186.900000000001
$ws.Range("N138").Value = -21280.0325

$ws.Range("H141").Value = 6231.3335
$ws.Range("I141").Value = 4183.75
$ws.Range("K141").Value = 12551.25
$ws.Range("M141").Value = -7371.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4906.029
$ws.Range("I32").Value = 4312.516
$ws.Range("J32").Value = 9505.75
$ws.Range("K32").Value = 4312.516
$ws.Range("L32").Value = 9505.75
$ws.Range("M32").Value = -4025.516
$ws.Range("N32").Value = -10079.75

$ws.Range("H63").Value = 5233.4165
$ws.Range("J63").Value = 8500
$ws.Range("L63").Value = 8500
$ws.Range("N63").Value = -9872

$ws.Range("H66").Value = 5233.4165
$ws.Range("J66").Value = 8500
$ws.Range("L66").Value = 42500
$ws.Range("N66").Value = -49364

$ws.Range("H96").Value = 39672
$ws.Range("J96").Value = 39672
$ws.Range("L96").Value = 39672
$ws.Range("N96").Value = -45164

$ws.Range("H101").Value = 30602
$ws.Range("J101").Value = 30602
$ws.Range("L101").Value = 30602
$ws.Range("N101").Value = -37092

$ws.Range("H102").Value = 2012.5
$ws.Range("I102").Value = 1765.75
$ws.Range("J102").Value = 2999.5
$ws.Range("K102").Value = 1765.75
$ws.Range("L102").Value = 2999.5
$ws.Range("M102").Value = -143.75
$ws.Range("N102").Value = -6243.5

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 33162.375
$ws.Range("I105").Value = 35883.332
$ws.Range("J105").Value = 24999.5
$ws.Range("K105").Value = 35883.332
$ws.Range("L105").Value = 24999.5
$ws.Range("M105").Value = -34136.332
$ws.Range("N105").Value = -28493.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25397.285
$ws.Range("I31").Value = 2313.4194
$ws.Range("J31").Value = 65152.832
$ws.Range("K31").Value = 2313.4194
$ws.Range("L31").Value = 65152.832
$ws.Range("M31").Value = -2018.4194
$ws.Range("N31").Value = -65742.83199999999

$ws.Range("H34").Value = 25397.285
$ws.Range("I34").Value = 2313.4194
$ws.Range("J34").Value = 65152.832
$ws.Range("K34").Value = 2313.4194
$ws.Range("L34").Value = 65152.832
$ws.Range("M34").Value = -2111.4194
$ws.Range("N34").Value = -65556.83199999999

$ws.Range("H62").Value = 5104
$ws.Range("I62").Value = 2004.25
$ws.Range("J62").Value = 17503
$ws.Range("K62").Value = 2004.25
$ws.Range("L62").Value = 17503
$ws.Range("M62").Value = -1380.25
$ws.Range("N62").Value = -18751

$ws.Range("H65").Value = 5104
$ws.Range("I65").Value = 2004.25
$ws.Range("J65").Value = 17503
$ws.Range("K65").Value = 10021.25
$ws.Range("L65").Value = 87515
$ws.Range("M65").Value = -6901.25
$ws.Range("N65").Value = -93755

$ws.Range("H86").Value = 11066.909
$ws.Range("I86").Value = 7549.8
$ws.Range("K86").Value = 7549.8
$ws.Range("M86").Value = -6426.8

$ws.Range("H89").Value = 11066.909
$ws.Range("I89").Value = 7549.8
$ws.Range("K89").Value = 37749
$ws.Range("M89").Value = -32133

$ws.Range("H99").Value = 3045
$ws.Range("I99").Value = 2590
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 2590
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -1092
$ws.Range("N99").Value = -6496

$ws.Range("H102").Value = 49999
$ws.Range("J102").Value = 49999
$ws.Range("L102").Value = 49999
$ws.Range("N102").Value = -54867

$ws.Range("H126").Value = 3045
$ws.Range("I126").Value = 2590
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 7770
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -5300
$ws.Range("N126").Value = -15440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 66760.53
$ws.Range("I6").Value = 83425.75
$ws.Range("J6").Value = 99.666664
$ws.Range("K6").Value = 250277.25
$ws.Range("L6").Value = 298.999992
$ws.Range("M6").Value = -250164.25
$ws.Range("N6").Value = -524.999992

$ws.Range("H39").Value = 3540
$ws.Range("I39").Value = 4000
$ws.Range("J39").Value = 3425
$ws.Range("K39").Value = 12000
$ws.Range("L39").Value = 10275
$ws.Range("M39").Value = -11706
$ws.Range("N39").Value = -10863

$ws.Range("H68").Value = 2974.9092
$ws.Range("J68").Value = 3452.889
$ws.Range("L68").Value = 10358.667
$ws.Range("N68").Value = -11980.667

$ws.Range("H71").Value = 2974.9092
$ws.Range("J71").Value = 3452.889
$ws.Range("L71").Value = 31076.001
$ws.Range("N71").Value = -39188.001

$ws.Range("H113").Value = 1350.1111
$ws.Range("I113").Value = 1616
$ws.Range("K113").Value = 4848
$ws.Range("M113").Value = -2678

$ws.Range("H125").Value = 8728.429
$ws.Range("I125").Value = 4000
$ws.Range("J125").Value = 9516.5
$ws.Range("K125").Value = 12000
$ws.Range("L125").Value = 28549.5
$ws.Range("M125").Value = -7080
$ws.Range("N125").Value = -38389.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

$ws.Range("H46").Value = 5422.923
$ws.Range("J46").Value = 5583.1665
$ws.Range("L46").Value = 5583.1665
$ws.Range("N46").Value = -5959.1665

$ws.Range("H55").Value = 2500776
$ws.Range("I55").Value = 3846813.8
$ws.Range("K55").Value = 3846813.8
$ws.Range("M55").Value = -3846640.8

$ws.Range("H82").Value = 5862.3335
$ws.Range("I82").Value = 1921
$ws.Range("K82").Value = 1921
$ws.Range("M82").Value = -1560

$ws.Range("H85").Value = 5862.3335
$ws.Range("I85").Value = 1921
$ws.Range("K85").Value = 1921
$ws.Range("M85").Value = -673

$ws.Range("H100").Value = 3233.7856
$ws.Range("I100").Value = 2168.4167
$ws.Range("K100").Value = 2168.4167
$ws.Range("M100").Value = -1627.4167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 35000
$ws.Range("J82").Value = 35000
$ws.Range("L82").Value = 35000
$ws.Range("N82").Value = -35766

$ws.Range("H85").Value = 35000
$ws.Range("J85").Value = 35000
$ws.Range("L85").Value = 35000
$ws.Range("N85").Value = -37652

$ws.Range("H113").Value = 531.55554
$ws.Range("I113").Value = 473
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1419
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 751
$ws.Range("N113").Value = -7340

$ws.Range("H135").Value = 59302.777
$ws.Range("J135").Value = 59302.777
$ws.Range("L135").Value = 59302.777
$ws.Range("N135").Value = -69442.777
